# Insert a new weekly record at row 152 (Vega Central Mapocho de Santiago - Rabanito),
# shifting the existing rows 152:177 down to 153:178.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data (old rows 152-177) down by one row.
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row 152 with the new weekly observation.
$ws.Range("A152").Value = 9
$ws.Range("B152").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C152").Value = "Metropolitana"
$ws.Range("D152").Value = 44504
$ws.Range("E152").Value = 13
$ws.Range("F152").Value = 300000001
$ws.Range("G152").Value = "Rabanito"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 7900
$ws.Range("K152").Value = 2500
$ws.Range("L152").Value = 3000
$ws.Range("M152").Value = 2747
$ws.Range("N152").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O152").Value = "Provincia de Chacabuco"
$ws.Range("P152").Value = 27
$ws.Range("Q152").Value = 100
$ws.Range("R152").Value = "Hortaliza"
